$d = $word.ActiveDocument

# 1) Split the paragraph: insert a paragraph break right after "1.2. " (and before
#    "Наименование работы") so the original run's tail text moves into a new paragraph,
#    without disturbing the later runs in that paragraph.
$splitPoint = $d.Content
$splitPoint.Find.Execute(
    "1.2. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$splitPoint.Collapse(0)
$splitPoint.InsertParagraphAfter()

# 2) Update the cached page-number field text in the header from "2" to "1".
$hdr = $d.Sections(1).Headers(1)
foreach ($f in $hdr.Range.Fields) {
    if ($f.Result.Text -eq "2") {
        $f.Result.Text = "1"
    }
}
